# Team_08_API Architects LMSTestData.xlsx - "Batch" sheet updates
# Adds an "updateBatchStatus" column (G) with an "Inactive" sample value,
# and updates the sheet's current selection to G3 (matching the authored
# commit: "schema validations and other batch validations").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Batch")
$ws.Activate()

# New header cell for the "updateBatchStatus" field (keeps the existing
# header styling already applied to G1).
$ws.Range("G1").Value = "updateBatchStatus"

# Sample/test data row for the new field.
$ws.Range("G2").Value = "Inactive"

# Move the sheet's active selection to G3, as captured in the saved view.
$ws.Range("G3").Select()
